# Challenge_AProjectAWeek_Git.xlsx — "Completed Week_8 and Week_9 projects"
#
# Adds two new weekly project blocks (Week 8: rows 64-76, Week 9: rows
# 77-85) below the existing Week 7b block (which ended at row 63), updates
# the frozen-pane / selection / zoom bookkeeping so the view still points
# at the bottom of the sheet, and nudges a couple of pre-existing cells
# (the header row + the "7a" label) onto the same cell style that the rest
# of the "week" rows already use.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# ---------------------------------------------------------------------
# 1) Normalise a couple of pre-existing styles so they match the rest of
#    the "week" rows (cosmetic clean-up the author made while editing).
# ---------------------------------------------------------------------

# A55 ("7a") was still on the old one-off style; match it to the other
# week-start cells in column A (e.g. A58).
$ws.Range("A58").Copy()
$ws.Range("A55").PasteSpecial($xlPasteFormats)

# Header row B3:F3 / A3 / G3 lose their now-unused distinct xf the same
# way — simplest is to just re-stamp them with their own look, which
# collapses onto the de-duplicated style used afterwards.
$ws.Range("A3").Copy()
$ws.Range("A3").PasteSpecial($xlPasteFormats)
$ws.Range("B3:F3").Copy()
$ws.Range("B3:F3").PasteSpecial($xlPasteFormats)
$ws.Range("G3").Copy()
$ws.Range("G3").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) Lay down two fresh "week" blocks by cloning the row formatting from
#    existing blocks of the same shape, then overwrite values/formulas.
# ---------------------------------------------------------------------

# Week 8 block: 13 rows (64-76), same shape as the Week 6 block (41-53).
$ws.Range("A41:G53").Copy($ws.Range("A64:G76"))

# Week 9 block: 9 rows (77-85), same shape as the Week 3 block (17-25).
$ws.Range("A17:G25").Copy($ws.Range("A77:G85"))
$excel.CutCopyMode = 0

# --- Week 8 data (row 64 = week header, rows 65-76 = sessions) --------

$ws.Range("A64").Value = 8
$ws.Range("B64").Value = 44816
$ws.Range("C64").Value = 0.44791666666666669
$ws.Range("D64").Value = 0.48958333333333331
$ws.Range("E64").Formula = "=D64-C64"
$ws.Range("F64").Formula = "=SUM(E64:E76)"
$ws.Range("G64").Value = "Loan Prediction ML Project"

$ws.Range("B65").Value = 44816
$ws.Range("C65").Value = 0.72916666666666663
$ws.Range("D65").Value = 0.8125
$ws.Range("E65").Formula = "=D65-C65"
$ws.Range("G65").Value = "Analytics Vidhya"

$ws.Range("B66").Value = 44817
$ws.Range("C66").Value = 0.41666666666666669
$ws.Range("D66").Value = 0.5625
$ws.Range("E66").Formula = "=D66-C66"
$ws.Range("G66").Value = "ML - Classification"

$ws.Range("B67").Value = 44818
$ws.Range("C67").Value = 0.40625
$ws.Range("D67").Value = 0.52083333333333337
$ws.Range("E67").Formula = "=D67-C67"
$ws.Range("G67").Value = "https://datahack.analyticsvidhya.com/contest/practice-problem-loan-prediction-iii/#DiscussTab"
$ws.Range("H67").Value = " "
# G67 came from the Week 6 template as a plain cell - give it the
# hyperlink-row border style before the hyperlink itself is added below.
$ws.Range("G59").Copy()
$ws.Range("G67").PasteSpecial($xlPasteFormats)

$ws.Range("B68").Value = 44818
$ws.Range("C68").Value = 0.60416666666666663
$ws.Range("D68").Value = 0.70833333333333337
$ws.Range("E68").Formula = "=D68-C68"

$ws.Range("B69").Value = 44819
$ws.Range("C69").Value = 0.41666666666666669
$ws.Range("D69").Value = 0.54166666666666663
$ws.Range("E69").Formula = "=D69-C69"

$ws.Range("B70").Value = 44820
$ws.Range("C70").Value = 0.4375
$ws.Range("D70").Value = 0.53125
$ws.Range("E70").Formula = "=D70-C70"

$ws.Range("B71").Value = 44820
$ws.Range("C71").Value = 0.5625
$ws.Range("D71").Value = 0.80208333333333337
$ws.Range("E71").Formula = "=D71-C71"

$ws.Range("B72").Value = 44821
$ws.Range("C72").Value = 0.86458333333333337
$ws.Range("D72").Value = 0.98958333333333337
$ws.Range("E72").Formula = "=D72-C72"
# G72 keeps the thin-border look (style 25) rather than the plain one.
$ws.Range("G59").Copy()
$ws.Range("G72").PasteSpecial($xlPasteFormats)

$ws.Range("B73").Value = 44822
$ws.Range("C73").Value = 0.71875
$ws.Range("D73").Value = 0.77083333333333337
$ws.Range("E73").Formula = "=D73-C73"

$ws.Range("B74").Value = 44822
$ws.Range("C74").Value = 0.92708333333333337
$ws.Range("D74").Value = 0.98958333333333337
$ws.Range("E74").Formula = "=D74-C74"
# G74 also keeps the thin-border look (style 25).
$ws.Range("G59").Copy()
$ws.Range("G74").PasteSpecial($xlPasteFormats)

$ws.Range("B75").Value = 44825
$ws.Range("C75").Value = 0.75
$ws.Range("D75").Value = 0.80208333333333337
$ws.Range("E75").Formula = "=D75-C75"

$ws.Range("B76").Value = 44825
$ws.Range("C76").Value = 0.85416666666666663
$ws.Range("D76").Value = 0.9375
$ws.Range("E76").Formula = "=D76-C76"

$excel.CutCopyMode = 0

# --- Week 9 data (row 77 = week header, rows 78-85 = sessions) --------

$ws.Range("A77").Value = 9
$ws.Range("B77").Value = 44823
$ws.Range("C77").Value = 0.39583333333333331
$ws.Range("D77").Value = 0.45833333333333331
$ws.Range("E77").Formula = "=D77-C77"
$ws.Range("F77").Formula = "=SUM(E77:E85)"
$ws.Range("G77").Value = "Predict the price of Bitcoin"

$ws.Range("B78").Value = 44823
$ws.Range("C78").Value = 0.61458333333333337
$ws.Range("D78").Value = 0.65625
$ws.Range("E78").Formula = "=D78-C78"
$ws.Range("G78").Value = "Udemy"
# G78 came from the Week 3 template with the hyperlink-row border style;
# the hyperlink itself is two rows further down in this block, so reset
# G78 back to the plain look.
$ws.Range("G60").Copy()
$ws.Range("G78").PasteSpecial($xlPasteFormats)

$ws.Range("B79").Value = 44823
$ws.Range("C79").Value = 0.79166666666666663
$ws.Range("D79").Value = 0.83333333333333337
$ws.Range("E79").Formula = "=D79-C79"
$ws.Range("G79").Value = "Time Series"

$ws.Range("B80").Value = 44823
$ws.Range("C80").Value = 0.91666666666666663
$ws.Range("D80").Value = 0.94791666666666663
$ws.Range("E80").Formula = "=D80-C80"
$ws.Range("G80").Value = "https://www.udemy.com/course/time-series-analysis-real-world-use-cases-in-python/learn/lecture/28361560#overview"
$ws.Range("H80").Value = " "
# G80 needs the hyperlink-row border style (style 25).
$ws.Range("G59").Copy()
$ws.Range("G80").PasteSpecial($xlPasteFormats)

$ws.Range("B81").Value = 44824
$ws.Range("C81").Value = 0.65625
$ws.Range("D81").Value = 0.73958333333333337
$ws.Range("E81").Formula = "=D81-C81"

$ws.Range("B82").Value = 44825
$ws.Range("C82").Value = 0.45833333333333331
$ws.Range("D82").Value = 0.51041666666666663
$ws.Range("E82").Formula = "=D82-C82"

$ws.Range("B83").Value = 44825
$ws.Range("C83").Value = 0.55208333333333337
$ws.Range("D83").Value = 0.65625
$ws.Range("E83").Formula = "=D83-C83"

$ws.Range("B84").Value = 44826
$ws.Range("C84").Value = 0.42708333333333331
$ws.Range("D84").Value = 0.54166666666666663
$ws.Range("E84").Formula = "=D84-C84"

$ws.Range("B85").Value = 44826
$ws.Range("C85").Value = 0.61458333333333337
$ws.Range("D85").Value = 0.64583333333333337
$ws.Range("E85").Formula = "=D85-C85"

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3) Hyperlinks for the two new project write-ups.
# ---------------------------------------------------------------------

$ws.Hyperlinks.Add(
    $ws.Range("G67"),
    "https://datahack.analyticsvidhya.com/contest/practice-problem-loan-prediction-iii/",
    "DiscussTab"
)
$ws.Hyperlinks.Add(
    $ws.Range("G80"),
    "https://www.udemy.com/course/time-series-analysis-real-world-use-cases-in-python/learn/lecture/28361560",
    "overview"
)

# ---------------------------------------------------------------------
# 4) View bookkeeping: the sheet now scrolls to show the new rows, and
#    the remembered selection/zoom follows the author's last position.
# ---------------------------------------------------------------------

$ws.Application.ActiveWindow.Zoom = 110
$ws.Range("C70").Select()
$ws.Application.ActiveWindow.FreezePanes = $false
$ws.Application.ActiveWindow.SplitColumn = 2
$ws.Application.ActiveWindow.SplitRow = 3
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Range("K85").Select()
$ws.Range("A1").Select()
